# Saldo_guide.xlsx update: roll the reference date forward one day
# (2024-07-30 -> 2024-07-31) and refresh the handful of balances that
# shifted as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new extraction timestamp.
$ws.Name = "IClientBalance-20240731-090451-"

# Column G ("Dt. Referencia") holds the reference date as a serial number
# for every data row (2-274). Excel serial 45503 = 2024-07-30,
# 45504 = 2024-07-31 -- bump every one of them.
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45504
}

# A few rows also had their projected/expected balances corrected; update
# columns D (Vl. Projetado), E (Saldo Previsto) and H (Vl. Total) to match.
$ws.Cells.Item(102, 5).Value = 15726.16   # E102
$ws.Cells.Item(102, 8).Value = 15726.16   # H102

$ws.Cells.Item(104, 4).Value = 0          # D104
$ws.Cells.Item(104, 5).Value = 33915.79   # E104
$ws.Cells.Item(104, 8).Value = 33915.79   # H104

$ws.Cells.Item(108, 4).Value = 0          # D108
$ws.Cells.Item(108, 5).Value = -16329.12  # E108

$ws.Cells.Item(173, 4).Value = 0          # D173
$ws.Cells.Item(173, 5).Value = 13061.83   # E173

$ws.Cells.Item(235, 5).Value = 365.53     # E235
$ws.Cells.Item(235, 8).Value = 365.53     # H235

$ws.Cells.Item(264, 4).Value = 0          # D264
$ws.Cells.Item(264, 5).Value = 1060.06    # E264
